$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "Example Aggregate Stats"

# Establish shared-string insertion order to match target workbook exactly.
$ws3.Range("A1").Value = "player_id"
$ws3.Range("B1").Value = "previous_player_id"
$ws3.Range("C1").Value = "live"

$ws3.Range("A2").Value = 1
$ws3.Range("B2").Value = "nil"
$ws3.Range("C2").Value = $false

$ws3.Range("E1").Value = "Player 1 Live Cells"
$ws3.Range("F1").Value = "Player 1 Dead Cells"
$ws3.Range("G1").Value = "Player 1 Regenerated Cells"
$ws3.Range("I1").Value = "Player 2 Dead Cells"
$ws3.Range("J1").Value = "Player 2 Regenerated Cells"
$ws3.Range("L1").Value = "Comments"

$ws3.Range("E2").Value = 0
$ws3.Range("F2").Value = 0
$ws3.Range("G2").Value = 0
$ws3.Range("H2").Value = 0
$ws3.Range("I2").Value = 0
$ws3.Range("J2").Value = 0
$ws3.Range("L2").Value = "Should not be possible. If the cell dies then the player id should move to previous player id"

$ws3.Range("A3").Value = 1
$ws3.Range("B3").Value = "nil"
$ws3.Range("C3").Value = $true
$ws3.Range("E3").Value = 1
$ws3.Range("F3").Value = 0
$ws3.Range("G3").Value = 0
$ws3.Range("H3").Value = 0
$ws3.Range("I3").Value = 0
$ws3.Range("J3").Value = 0

$ws3.Range("A4").Value = "nil"
$ws3.Range("B4").Value = 1
$ws3.Range("C4").Value = $false
$ws3.Range("E4").Value = 0
$ws3.Range("F4").Value = 1
$ws3.Range("G4").Value = 0
$ws3.Range("H4").Value = 0
$ws3.Range("I4").Value = 0
$ws3.Range("J4").Value = 0

$ws3.Range("A5").Value = "nil"
$ws3.Range("B5").Value = 1
$ws3.Range("C5").Value = $true
$ws3.Range("E5").Value = 0
$ws3.Range("F5").Value = 0
$ws3.Range("G5").Value = 0
$ws3.Range("H5").Value = 0
$ws3.Range("I5").Value = 0
$ws3.Range("J5").Value = 0
$ws3.Range("L5").Value = "Should not be possible. If the cell is alive then the player_id should be set"

$ws3.Range("A6").Value = 1
$ws3.Range("B6").Value = 1
$ws3.Range("C6").Value = $false
$ws3.Range("E6").Value = 0
$ws3.Range("F6").Value = 0
$ws3.Range("G6").Value = 0
$ws3.Range("H6").Value = 0
$ws3.Range("I6").Value = 0
$ws3.Range("J6").Value = 0

$ws3.Range("A7").Value = 1
$ws3.Range("B7").Value = 1
$ws3.Range("C7").Value = $true
$ws3.Range("E7").Value = 1
$ws3.Range("F7").Value = 0
$ws3.Range("G7").Value = 1
$ws3.Range("H7").Value = 0
$ws3.Range("I7").Value = 0
$ws3.Range("J7").Value = 0

$ws3.Range("A8").Value = 1
$ws3.Range("B8").Value = 2
$ws3.Range("C8").Value = $false
$ws3.Range("E8").Value = 0
$ws3.Range("F8").Value = 0
$ws3.Range("G8").Value = 0
$ws3.Range("H8").Value = 0
$ws3.Range("I8").Value = 0
$ws3.Range("J8").Value = 0

# H1 header was entered after the row-2/5 comments but before the row-8
# comment (matches shared-string table order in the target file).
$ws3.Range("H1").Value = "Player 2 Live Cells"

$ws3.Range("L8").Value = "Should not be possible. If the cell is dead then there should be no player_id set"

$ws3.Range("A9").Value = 1
$ws3.Range("B9").Value = 2
$ws3.Range("C9").Value = $true
$ws3.Range("E9").Value = 1
$ws3.Range("F9").Value = 0
$ws3.Range("G9").Value = 1
$ws3.Range("H9").Value = 0
$ws3.Range("I9").Value = 0
$ws3.Range("J9").Value = 0

$ws3.Columns.AutoFit()
$ws3.Range("F21").Select()
$ws3.Activate()
